$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# The "fecha" (date) column had previously been inserted as column A by an
# export helper that stamped every row with the export date. This commit
# removes that auto-added "fecha" column again (column A), so the sheet
# goes back to Id/Nombre/Materia/Turno/.../Ins_Cupo starting at column B.
$ws.Columns.Item(1).Clear()

# Keep the last selected cell as recorded by the author when saving.
$ws.Range("C14").Select()
